$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "_old" header columns (A:J) to "_FV2210" and the
# "_new" header columns (L:U) to "_FV2304". Column K ("diff") is untouched.
$headerBases = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $headerBases.Length; $i++) {
    $oldCol = $i + 1        # columns 1..10 -> A..J
    $newCol = $i + 12       # columns 12..21 -> L..U
    $ws.Cells.Item(1, $oldCol).Value = "$($headerBases[$i])_FV2210"
    $ws.Cells.Item(1, $newCol).Value = "$($headerBases[$i])_FV2304"
}

# Turn the data range into an Excel Table ("Table1") with headers.
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U52"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# Freeze the header row (split below row 1, freeze top pane).
$ws.Range("A2").Select() | Out-Null
($excel.ActiveWindow.FreezePanes = $true) | Out-Null
